$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B34").Value = "9b5fa738b68a8c46f512c3e8ae609d3b"
$ws.Range("B44").Value = "a2cfcbfef9b7b4aed5ed06cdf76e820f"
$ws.Range("B74").Value = "9555bf74da8a390313ded720eb47dce7"
$ws.Range("B89").Value = "160ee88f449d69ffbf488ebe9d2dcc44"
$ws.Range("B99").Value = "ec5bd2a050b8a245967e920be6cdaaa2"
$ws.Range("B110").Value = "4050bd447a74401c61ea746f9711d4fc"
$ws.Range("B154").Value = "0164192226833e8b2508d9634b0ba903"
$ws.Range("B160").Value = "adf3c1215f1ec05392a34e4fcab6d818"
$ws.Range("B161").Value = "9bb4c7968671c6ffbee5b3db18131f17"
$ws.Range("B162").Value = "537a5222143850acb0b8e7c2a56d1a6f"
$ws.Range("B168").Value = "36c8cd53ba8a46717318adc0a51706b1"
$ws.Range("B180").Value = "8e3e66726412138b9c21d57bc4009d98"
$ws.Range("B213").Value = "f1a3da6a4991d211f4d0e18b9486ed7a"
$ws.Range("B222").Value = "611f46935622466783f30ce5419ec305"
$ws.Range("B229").Value = "b946e436d07d4b85b5db0149d8365a58"
$ws.Range("B278").Value = "4f4e6e1d7f91885a3a4f184b8ac396e3"
$ws.Range("B330").Value = "0f541db1bee54323ba63ade30ce40dfc"
$ws.Range("B345").Value = "183913fecc02620ae6913e0667b17656"
$ws.Range("B461").Value = "060072cb4a449d58d07838c00b609f70"
$ws.Range("B506").Value = "aa1791820592e49d2dde3aff5748084a"
$ws.Range("B514").Value = "0163ad4ebad868ebcb1fb1d515410e6b"
$ws.Range("B516").Value = "c239325c8f39fbe763a8255fd7fd1e8d"
$ws.Range("B524").Value = "b8463e643f40c14c051b7aa3e19cc647"
$ws.Range("B534").Value = "b4d216af1c0225064ccc574065e16246"
$ws.Range("B547").Value = "12134a6651c6de21c72dc6c1e1dae89a"
$ws.Range("B553").Value = "8317bc5e1079993b6d686cc7d773b4ef"
$ws.Range("B666").Value = "d0198b482e7ad0701fea272aba6657a8"
$ws.Range("B729").Value = "b4db0bd5cfe9f51ea71702c7935a8b82"
$ws.Range("B768").Value = "8a866f38cea4d509d812189b47eef642"
$ws.Range("B811").Value = "5f1e48ea2ee37ac4a0cd6534daf28e1d"
$ws.Range("B815").Value = "deeeabb02d47e448e34e5d3bbaeb8dad"
$ws.Range("B816").Value = "1951623ae9020a139ec3467817acc2ab"
$ws.Range("B825").Value = "76fb08e3968f1341beee8c4d704ab1a6"
$ws.Range("B827").Value = "fe391b223dd9b3e7fc6a5f6ebd9890a3"
